# Add the "fillForm" worksheet after "search" and populate it with the
# submitted form data, then select it as the active sheet (matches
# activeTab="1" / tabSelected flip in the target workbook).

$wb = $excel.ActiveWorkbook
$search = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add($null, $search)
$ws.Name = "fillForm"

# Whole column A is formatted as Text ("@") before any values are written so
# that every entry - including the numeric-looking zip code - is stored as
# text, matching the source data.
$ws.Range("A1:A10").NumberFormat = "@"

$ws.Range("A1").Value = "Form Parameters"
$ws.Range("A2").Value = "nellore.mahesh867@gmail.com"
$ws.Range("A3").Value = "Breville"
$ws.Range("A4").Value = "Test"
$ws.Range("A5").Value = "900 Hamlin"
$ws.Range("A6").Value = "Ct"
$ws.Range("A7").Value = "Sunnyvale"
$ws.Range("A8").Value = "94089"
$ws.Range("A9").Value = "California"
$ws.Range("A10").Value = "(541) 754-3010"

# Header row gets a cyan fill.
$ws.Range("A1").Interior.PatternColor = 13421619
$ws.Range("A1").Interior.Color = 16763904

# Email cell becomes a live mailto: hyperlink.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:nellore.mahesh867@gmail.com", "", "", "nellore.mahesh867@gmail.com")

# Zip code is left-aligned.
$ws.Range("A8").HorizontalAlignment = -4131

$ws.Columns.Item(1).ColumnWidth = 26.05

$ws.Range("B10").Select()
